$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1091
$ws1.Range("F5").Value = 3679
$ws1.Range("F7").Value = 304
$ws1.Range("F9").Value = 9
$ws1.Range("F12").Value = 112
$ws1.Range("F13").Value = 281
$ws1.Range("G13").Value = 49
$ws1.Range("F14").Value = 54
$ws1.Range("F15").Value = 106
$ws1.Range("F16").Value = 2820
$ws1.Range("F17").Value = 1153
$ws1.Range("F18").Value = 7

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1091
$ws4.Range("F6").Value = 3679
$ws4.Range("F8").Value = 304
$ws4.Range("F11").Value = 9
$ws4.Range("F14").Value = 112
$ws4.Range("F15").Value = 281
$ws4.Range("G15").Value = 49
$ws4.Range("F16").Value = 54
$ws4.Range("F17").Value = 106
$ws4.Range("F18").Value = 2820
$ws4.Range("F19").Value = 1153
$ws4.Range("F20").Value = 7
